# Remove the "controleren van de configuratie op aanwezigheid van bekende
# kwetsbaarheden," bullet paragraph from the tools list on slide 13.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(13)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

$target = "controleren van de configuratie op aanwezigheid van bekende kwetsbaarheden,"

for ($i = $tr.Paragraphs().Count; $i -ge 1; $i--) {
    $para = $tr.Paragraphs($i)
    if ($para.Text.Trim() -eq $target) {
        $para.Delete()
    }
}
